$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.094361782073975
$ws.Range("B1").Value = 1.892000555992126
$ws.Range("D1").Value = 1.187701225280762
$ws.Range("E1").Value = 1.157837748527527
